# Updated symbol list on Thu Jan  5 07:28:24 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns with the latest scrape,
# and fixes the swapped BKEXToken / KickToken rows (41-42).
#
# NOTE: values in D/E are plain text (e.g. "4.750", "1.24%") rather than
# numbers, so NumberFormat is forced to "@" (Text) before each assignment -
# otherwise COM auto-converts numeric-looking strings to real numbers and
# silently drops significant trailing zeros (e.g. "3.480" -> 3.48) or turns
# percentages into their decimal fraction (e.g. "1.04%" -> 0.0104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "257.31"
Set-TextValue "E2" "1.04%"

# Row 3 - OKB
Set-TextValue "D3" "26.97"
Set-TextValue "E3" "-3.79%"

# Row 4 - HuobiToken
Set-TextValue "D4" "4.728"
Set-TextValue "E4" "-11.29%"

# Row 5 - Cronos
Set-TextValue "D5" "0.05972"
Set-TextValue "E5" "2.10%"

# Row 6 - KuCoinToken
Set-TextValue "D6" "6.675"
Set-TextValue "E6" "-0.61%"

# Row 7 - MXToken
Set-TextValue "D7" "0.8693"
Set-TextValue "E7" "0.32%"

# Row 8 - FTXToken
Set-TextValue "D8" "0.9445"
Set-TextValue "E8" "3.72%"

# Row 9 - WazirX
Set-TextValue "D9" "0.1407"
Set-TextValue "E9" "-1.17%"

# Row 10 - LiechtensteinCryptoassetsExchange
Set-TextValue "D10" "0.03625"
Set-TextValue "E10" "4.33%"

# Row 11 - MandalaExchangeToken
Set-TextValue "D11" "0.07183"
Set-TextValue "E11" "0.31%"

# Row 12 - BitrueCoin
Set-TextValue "D12" "0.03176"
Set-TextValue "E12" "-0.16%"

# Row 13 - BitMartToken
Set-TextValue "D13" "0.09241"
Set-TextValue "E13" "0.23%"

# Row 14 - BitForexToken
Set-TextValue "D14" "0.001541"
Set-TextValue "E14" "0.09%"

# Row 15 - One
Set-TextValue "D15" "0.0006082"
Set-TextValue "E15" "-94.28%"

# Row 16 - TigerCash
Set-TextValue "D16" "0.006098"
Set-TextValue "E16" "2.95%"

# Row 17 - LEO
Set-TextValue "D17" "3.480"
Set-TextValue "E17" "-0.57%"

# Row 18 - GateToken
Set-TextValue "D18" "3.182"
Set-TextValue "E18" "-1.50%"

# Row 19 - BTSEToken
Set-TextValue "D19" "2.240"
Set-TextValue "E19" "1.72%"

# Row 20 - BitpandaEcosystemToken (price unchanged)
Set-TextValue "E20" "-1.95%"

# Row 21 - ProBitToken
Set-TextValue "D21" "0.1307"
Set-TextValue "E21" "-0.68%"

# Row 22 - MCDex
Set-TextValue "D22" "3.826"
Set-TextValue "E22" "7.47%"

# Row 23 - CoinExToken
Set-TextValue "D23" "0.04219"
Set-TextValue "E23" "1.53%"

# Row 25 - BitKan
Set-TextValue "D25" "0.001226"
Set-TextValue "E25" "0.19%"

# Row 26 - HotbitToken
Set-TextValue "D26" "0.004501"
Set-TextValue "E26" "-10.67%"

# Row 27 - NitroEx (price unchanged)
Set-TextValue "E27" "41.84%"

# Row 28 - UpBots (price unchanged)
Set-TextValue "E28" "-22.91%"

# Row 40 - IDEX
Set-TextValue "D40" "0.03817"
Set-TextValue "E40" "-0.78%"

# Row 41 - was BKEXToken, is now KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006165"
Set-TextValue "E41" "8.77%"

# Row 42 - was KickToken, is now BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1099"
Set-TextValue "E42" "-0.04%"

# Row 43 - CEJI
Set-TextValue "D43" "0.002254"
Set-TextValue "E43" "2.47%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.01059"
Set-TextValue "E44" "-2.98%"

# Row 45 - CoinLion (price unchanged)
Set-TextValue "E45" "5.07%"

# Row 46 - Kangarootoken (volume unchanged)
Set-TextValue "D46" "0.00000000751"

# Row 47 - CoinbaseStockToken (price unchanged)
Set-TextValue "E47" "21.50%"

# Row 48 - BOLO
Set-TextValue "D48" "0.002278"
Set-TextValue "E48" "5.64%"

# Row 49 - CryptobidCoin (volume unchanged)
Set-TextValue "D49" "0.00002102"

# Row 50 - SpecialPowerGold (volume unchanged)
Set-TextValue "D50" "0.0002002"
